$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("reaction8")

# Clear out the old row of values beyond column B
$ws.Range("C1:O1").ClearContents()

# Set the new values for the remaining cells
$ws.Range("A1").Value = 16
$ws.Range("B1").Value = 17
